# Weekly NYPD CompStat crime data update (cs-en-us-044pct.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper: write a text placeholder ("0" / "***.*") into a cell that was
# previously numeric, matching the General-format/text style (s=13) used by
# the other placeholder cells in the table (e.g. A14).
function Set-TextPlaceholder($addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value() = $text
    $ws.Range("A14").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Header: volume/issue number and reporting week dates
$ws.Range("A8").Value() = "Volume 32   Number  45"
$ws.Range("C9").Value() = "Report Covering the Week  11/3/2025  Through  11/9/2025"

# Crime statistics table (rows 14-33)
$ws.Range("L14").Value() = -46.666666666666
$ws.Range("N14").Value() = -86.885245901639
Set-TextPlaceholder "C15" "0"
$ws.Range("D15").Value() = 4
$ws.Range("E15").Value() = -100
$ws.Range("F15").Value() = 4
$ws.Range("G15").Value() = 6
$ws.Range("H15").Value() = -33.333333333333
$ws.Range("I15").Value() = 34
$ws.Range("J15").Value() = 38
$ws.Range("K15").Value() = -10.526315789473
$ws.Range("L15").Value() = -2.857142857142
$ws.Range("M15").Value() = 21.428571428571
$ws.Range("N15").Value() = -56.962025316455
$ws.Range("C16").Value() = 4
$ws.Range("E16").Value() = -71.428571428571
$ws.Range("F16").Value() = 32
$ws.Range("G16").Value() = 53
$ws.Range("H16").Value() = -39.622641509434
$ws.Range("I16").Value() = 403
$ws.Range("J16").Value() = 495
$ws.Range("K16").Value() = -18.585858585858
$ws.Range("L16").Value() = -15.336134453781
$ws.Range("M16").Value() = -5.176470588235
$ws.Range("N16").Value() = -76.825761932144
$ws.Range("C17").Value() = 22
$ws.Range("D17").Value() = 21
$ws.Range("E17").Value() = 4.761904761904
$ws.Range("F17").Value() = 71
$ws.Range("G17").Value() = 82
$ws.Range("H17").Value() = -13.414634146341
$ws.Range("I17").Value() = 895
$ws.Range("J17").Value() = 923
$ws.Range("K17").Value() = -3.033586132177
$ws.Range("L17").Value() = -6.282722513089
$ws.Range("M17").Value() = 86.458333333333
$ws.Range("N17").Value() = -25.354462051709
$ws.Range("C18").Value() = 3
$ws.Range("D18").Value() = 7
$ws.Range("E18").Value() = -57.142857142857
$ws.Range("F18").Value() = 17
$ws.Range("G18").Value() = 26
$ws.Range("H18").Value() = -34.615384615384
$ws.Range("I18").Value() = 203
$ws.Range("J18").Value() = 257
$ws.Range("K18").Value() = -21.011673151751
$ws.Range("L18").Value() = -48.477157360406
$ws.Range("M18").Value() = -11.739130434782
$ws.Range("N18").Value() = -88.709677419354
$ws.Range("C19").Value() = 17
$ws.Range("D19").Value() = 23
$ws.Range("E19").Value() = -26.086956521739
$ws.Range("F19").Value() = 61
$ws.Range("G19").Value() = 91
$ws.Range("H19").Value() = -32.967032967033
$ws.Range("I19").Value() = 649
$ws.Range("J19").Value() = 783
$ws.Range("K19").Value() = -17.113665389527
$ws.Range("L19").Value() = -7.285714285714
$ws.Range("M19").Value() = 97.264437689969
$ws.Range("N19").Value() = -7.943262411347
$ws.Range("C20").Value() = 2
$ws.Range("D20").Value() = 8
$ws.Range("E20").Value() = -75
$ws.Range("F20").Value() = 13
$ws.Range("G20").Value() = 28
$ws.Range("H20").Value() = -53.571428571428
$ws.Range("I20").Value() = 281
$ws.Range("J20").Value() = 315
$ws.Range("K20").Value() = -10.793650793650
$ws.Range("L20").Value() = -31.796116504854
$ws.Range("M20").Value() = 53.551912568306
$ws.Range("N20").Value() = -75.285839929639
$ws.Range("C21").Value() = 48
$ws.Range("D21").Value() = 77
$ws.Range("E21").Value() = -37.662337662337
$ws.Range("F21").Value() = 198
$ws.Range("G21").Value() = 286
$ws.Range("H21").Value() = -30.769230769230
$ws.Range("I21").Value() = 2473
$ws.Range("J21").Value() = 2824
$ws.Range("K21").Value() = -12.429178470255
$ws.Range("L21").Value() = -17.207900903917
$ws.Range("M21").Value() = 46.504739336492
$ws.Range("N21").Value() = -63.188448943137
Set-TextPlaceholder "C22" "0"
$ws.Range("F22").Value() = 1
$ws.Range("G22").Value() = 2
$ws.Range("H22").Value() = -50
$ws.Range("I22").Value() = 36
$ws.Range("J22").Value() = 44
$ws.Range("K22").Value() = -18.181818181818
$ws.Range("L22").Value() = -5.263157894736
$ws.Range("M22").Value() = -2.702702702702
Set-TextPlaceholder "C23" "0"
Set-TextPlaceholder "D23" "0"
Set-TextPlaceholder "E23" "***.*"
$ws.Range("F23").Value() = 4
$ws.Range("G23").Value() = 3
$ws.Range("H23").Value() = 33.333333333333
$ws.Range("I23").Value() = 38
$ws.Range("J23").Value() = 66
$ws.Range("K23").Value() = -42.424242424242
$ws.Range("L23").Value() = -44.927536231884
$ws.Range("M23").Value() = -5
$ws.Range("C24").Value() = 17
$ws.Range("D24").Value() = 42
$ws.Range("E24").Value() = -59.523809523809
$ws.Range("F24").Value() = 99
$ws.Range("G24").Value() = 140
$ws.Range("H24").Value() = -29.285714285714
$ws.Range("I24").Value() = 1545
$ws.Range("J24").Value() = 1672
$ws.Range("K24").Value() = -7.595693779904
$ws.Range("L24").Value() = -13.783482142857
$ws.Range("M24").Value() = 29.614093959731
$ws.Range("C25").Value() = 6
$ws.Range("D25").Value() = 26
$ws.Range("E25").Value() = -76.923076923076
$ws.Range("F25").Value() = 33
$ws.Range("G25").Value() = 58
$ws.Range("H25").Value() = -43.103448275862
$ws.Range("I25").Value() = 575
$ws.Range("J25").Value() = 764
$ws.Range("K25").Value() = -24.738219895288
$ws.Range("L25").Value() = -38.105489773950
$ws.Range("C26").Value() = 19
$ws.Range("D26").Value() = 28
$ws.Range("E26").Value() = -32.142857142857
$ws.Range("F26").Value() = 97
$ws.Range("G26").Value() = 113
$ws.Range("H26").Value() = -14.159292035398
$ws.Range("I26").Value() = 1123
$ws.Range("J26").Value() = 1248
$ws.Range("K26").Value() = -10.016025641025
$ws.Range("L26").Value() = 4.855275443510
$ws.Range("M26").Value() = 3.027522935779
Set-TextPlaceholder "C27" "0"
$ws.Range("D27").Value() = 4
$ws.Range("E27").Value() = -100
$ws.Range("G27").Value() = 7
$ws.Range("H27").Value() = -14.285714285714
$ws.Range("I27").Value() = 52
$ws.Range("J27").Value() = 58
$ws.Range("K27").Value() = -10.344827586206
$ws.Range("L27").Value() = -14.754098360655
$ws.Range("C28").Value() = 6
$ws.Range("D28").Value() = 1
$ws.Range("E28").Value() = 500
$ws.Range("G28").Value() = 16
$ws.Range("H28").Value() = -18.75
$ws.Range("I28").Value() = 117
$ws.Range("J28").Value() = 120
$ws.Range("K28").Value() = -2.5
$ws.Range("L28").Value() = 1.739130434782
Set-TextPlaceholder "C29" "0"
$ws.Range("F29").Value() = 1
$ws.Range("G29").Value() = 2
$ws.Range("H29").Value() = -50
$ws.Range("J29").Value() = 56
$ws.Range("K29").Value() = -44.642857142857
$ws.Range("M29").Value() = -16.216216216216
$ws.Range("N29").Value() = -80.745341614906
Set-TextPlaceholder "C30" "0"
$ws.Range("F30").Value() = 1
$ws.Range("G30").Value() = 2
$ws.Range("H30").Value() = -50
$ws.Range("J30").Value() = 41
$ws.Range("K30").Value() = -41.463414634146
$ws.Range("M30").Value() = -20
$ws.Range("N30").Value() = -83.448275862069
$ws.Range("L33").Value() = -60

$excel.CutCopyMode = 0
